$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.690954089164734
$ws.Range("B1").Value = 2.28130054473877
$ws.Range("C1").Value = 5.160366058349609
$ws.Range("D1").Value = 1.377427101135254
$ws.Range("E1").Value = 0.650273859500885
